$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.904.88'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.915.96'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.52'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.18'
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.85'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.438'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000225'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '33.54'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.398.33'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '60.890.01'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.68'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.916.42'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '429.58'
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.36'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.679'
$ws.Range('E21').Value = '  +1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.06'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '81.40'
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.91'
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.20'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.94'
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.29'
$ws.Range('E28').Value = '  +4.60%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.61'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.03'
$ws.Range('E31').Value = '  -2.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.36'
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0846'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.62'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.01'
$ws.Range('E37').Value = '  +1.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.98'
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.121'
$ws.Range('E39').Value = '  -2.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.51'
$ws.Range('E40').Value = '  -1.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.289'
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.87'
$ws.Range('E42').Value = '  -4.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '374.99'
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0344'
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.698.58'
$ws.Range('E45').Value = '  +1.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '131.20'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.87'
$ws.Range('E48').Value = '  -6.11%  '
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.01'
$ws.Range('E50').Value = '  -3.95%  '
$ws.Range('E51').Value = '  +1.48%  '
